$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data cells (row 7: Read/Study)
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2

# Row 8: Team Meting
$ws.Range("F8").Value = 1

# Row 11: Smaller Team Meet
$ws.Range("E11").Value = 0.5

# Update the active cell selection to H9
$ws.Range("H9").Select()
